$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHTR")

# Insert two new columns before column D, shifting existing data (D:K) to (F:M)
$ws.Range("D5:E5").EntireColumn.Insert()

# Copy cell formatting from the (now-shifted) original D:E columns (now F:G) into the new D:E columns,
# but only for the rows that actually held data in the three financial tables (rows 5,6,37,79 are bare
# section-label rows with no D:K cells at all, so they must stay untouched).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Populate the new D (most recent quarter) and E (prior quarter) columns with new data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 11231000
$ws.Range("E8").Value = 10892000
$ws.Range("D9").Value = 6407000
$ws.Range("E9").Value = 6222000
$ws.Range("D10").Value = 4824000
$ws.Range("E10").Value = 4670000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 25000
$ws.Range("E14").Value = 14000
$ws.Range("D15").Value = 2534000
$ws.Range("E15").Value = 2482000
$ws.Range("D17").Value = 9792000
$ws.Range("E17").Value = 9512000
$ws.Range("D18").Value = 1439000
$ws.Range("E18").Value = 1380000
$ws.Range("D20").Value = -167000
$ws.Range("E20").Value = 214000
$ws.Range("D21").Value = 3806000
$ws.Range("E21").Value = 4076000
$ws.Range("D22").Value = 910000
$ws.Range("E22").Value = 901000
$ws.Range("D23").Value = 362000
$ws.Range("E23").Value = 693000
$ws.Range("D24").Value = 2000
$ws.Range("E24").Value = 109000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 360000
$ws.Range("E26").Value = 584000
$ws.Range("D27").Value = 296000
$ws.Range("E27").Value = 493000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 167000
$ws.Range("E32").Value = -214000
$ws.Range("D33").Value = 296000
$ws.Range("E33").Value = 493000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 296000
$ws.Range("E35").Value = 493000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 551000
$ws.Range("E41").Value = 612000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1733000
$ws.Range("E43").Value = 1736000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 446000
$ws.Range("E45").Value = 381000
$ws.Range("D46").Value = 2730000
$ws.Range("E46").Value = 2729000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 35126000
$ws.Range("E48").Value = 34740000
$ws.Range("D49").Value = 106914000
$ws.Range("E49").Value = 107483000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1360000
$ws.Range("E52").Value = 1133000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 146130000
$ws.Range("E54").Value = 146085000
$ws.Range("D57").Value = 758000
$ws.Range("E57").Value = 604000
$ws.Range("D58").Value = 3290000
$ws.Range("E58").Value = 3339000
$ws.Range("D59").Value = 8047000
$ws.Range("E59").Value = 7907000
$ws.Range("D60").Value = 12095000
$ws.Range("E60").Value = 11850000
$ws.Range("D61").Value = 69537000
$ws.Range("E61").Value = 69135000
$ws.Range("D62").Value = 20226000
$ws.Range("E62").Value = 19872000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 109845000
$ws.Range("E66").Value = 108980000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2780000
$ws.Range("E72").Value = 4828000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 36285000
$ws.Range("E76").Value = 37105000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 296000
$ws.Range("E81").Value = 493000
$ws.Range("D83").Value = 2534000
$ws.Range("E83").Value = 2482000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 3168000
$ws.Range("E89").Value = 2804000
$ws.Range("D91").Value = -2433000
$ws.Range("E91").Value = -2118000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -2306000
$ws.Range("E94").Value = -2323000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -757000
$ws.Range("E100").Value = -594000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 105000
$ws.Range("E102").Value = -113000
